# ---------------------------------------------------------------------------
# Adds a "2022-Q4" quarterly sheet to the workbook (inserted between "总计"
# and "2022-Q3"), populates it with the fund-holding detail rows, and
# updates the "总计" (totals) sheet so the new quarter is reflected there.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1) Update the "总计" (totals) summary sheet (sheet1): a new 2022-Q4 row
#    is inserted at the top of the data table and every other quarter
#    shifts down by one row.
# -----------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$dates  = @("2022-Q4", "2022-Q3", "2022-Q2", "2022-Q1", "2021-Q4", "2021-Q3", "2021-Q2")
$counts = @(20, 4, 11, 8, 10, 7, 1)
$values = @(0.53, 0.09, 0.33, 0.25, 1.06, 0.24, 0.02)

# Give the new row 8 the same look (bold/bordered index cell) as the
# existing data rows before writing into it.
$total.Range("A7").Copy($total.Range("A8"))

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $i + 2
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $dates[$i]
    $total.Cells.Item($r, 3).Value = $counts[$i]
    $total.Cells.Item($r, 4).Value = $values[$i]
}

# -----------------------------------------------------------------------
# 2) Insert the new "2022-Q4" sheet right after "总计" (i.e. before the
#    current "2022-Q3" sheet) and populate it with fund-holding detail.
#    Duplicating the "2022-Q3" sheet (rather than Worksheets.Add) is what
#    carries its styling (bold/bordered header row + index column) over
#    to the new sheet, and Excel places the duplicate immediately before
#    the sheet it was copied from - exactly the slot "2022-Q4" needs.
# -----------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Extend the styled index column down to row 21 (the source sheet only
# had 5 data rows / row 21 is needed for the 20 funds below).
$q4.Range("A2").Copy($q4.Range("A2:A21"))

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$rows = @(
    @("010064", "圆信永丰兴研灵活配置混合A", "9.31", "88.84", "1.66", "0.1545", 10),
    @("001468", "广发改革先锋灵活配置混合", "5.72", "93.29", "2.39", "0.1367", 9),
    @("009847", "圆信永丰研究精选混合A", "1.13", "93.72", "4.73", "0.0534", 3),
    @("006969", "圆信永丰高端制造混合", "0.89", "92.44", "5.47", "0.0487", 2),
    @("009056", "圆信永丰大湾区主题混合C", "0.67", "90.09", "3.27", "0.0219", 9),
    @("004657", "金鹰民富收益混合A", "2.07", "37.03", "0.87", "0.0180", 5),
    @("010065", "圆信永丰兴研灵活配置混合C", "0.98", "88.84", "1.66", "0.0163", 10),
    @("010740", "汇安核心价值混合A", "0.41", "92.83", "3.74", "0.0153", 8),
    @("012498", "汇添富中证500基本面增强指数A", "0.77", "92.29", "1.51", "0.0116", 10),
    @("009055", "圆信永丰大湾区主题混合A", "0.34", "90.09", "3.27", "0.0111", 9),
    @("004658", "金鹰民富收益混合C", "1.18", "37.03", "0.87", "0.0103", 5),
    @("011433", "中加聚优一年定期开放混合A", "0.87", "24.92", "1.15", "0.0100", 7),
    @("010741", "汇安核心价值混合C", "0.20", "92.83", "3.74", "0.0075", 8),
    @("009054", "圆信永丰沣泰混合", "0.23", "31.90", "1.64", "0.0038", 6),
    @("012499", "汇添富中证500基本面增强指数C", "0.23", "92.29", "1.51", "0.0035", 10),
    @("009848", "圆信永丰研究精选混合C", "0.05", "93.72", "4.73", "0.0024", 3),
    @("008838", "德邦量化对冲策略灵活配置混合A", "0.16", "65.09", "0.95", "0.0015", 5),
    @("515590", "前海开源中证500等权重ETF", "0.34", "95.19", "0.29", "0.0010", 9),
    @("008839", "德邦量化对冲策略灵活配置混合C", "0.04", "65.09", "0.95", "0.0004", 5),
    @("011434", "中加聚优一年定期开放混合C", "0.03", "24.92", "1.15", "0.0003", 7)
)

# Fund codes (column B) and the numeric-looking ratio columns (D:G) are
# stored as plain text in the source workbook (e.g. "010064" keeps its
# leading zero) - force text format before writing so Excel doesn't
# auto-convert them to numbers, then clear the (otherwise unused) number
# format back off so the cells end up with the workbook's default style,
# exactly like the other quarterly sheets.
$q4.Range("B2:B21").NumberFormat = "@"
$q4.Range("D2:G21").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $q4.Cells.Item($r, 1).Value = $i
    $q4.Cells.Item($r, 2).Value = $row[0]
    $q4.Cells.Item($r, 3).Value = $row[1]
    $q4.Cells.Item($r, 4).Value = $row[2]
    $q4.Cells.Item($r, 5).Value = $row[3]
    $q4.Cells.Item($r, 6).Value = $row[4]
    $q4.Cells.Item($r, 7).Value = $row[5]
    $q4.Cells.Item($r, 8).Value = $row[6]
}

$q4.Range("B2:B21").ClearFormats()
$q4.Range("D2:G21").ClearFormats()

Write-Output "done"
